$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.114.91'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').Value = '2.281.36'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''155.52'
$ws.Range('E5').Value = '  +15,433.66%  '
$ws.Range('D6').Value = '''305.17'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').Value = '''94.81'
$ws.Range('E7').Value = '  +2.05%  '
$ws.Range('E8').Value = '  -0.25%  '
$ws.Range('E10').Value = '  +1.15%  '
$ws.Range('D11').Value = '''35.34'
$ws.Range('E11').Value = '  +8.13%  '
$ws.Range('D12').Value = '''0.0804'
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('E13').Value = '  -2.01%  '
$ws.Range('D15').Value = '2.634.81'
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('D16').Value = '''14.43'
$ws.Range('E16').Value = '  +1.23%  '
$ws.Range('D17').Value = '2.275.24'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('E18').Value = '  +4.23%  '
$ws.Range('D19').Value = '42.056.53'
$ws.Range('D20').Value = '''12.81'
$ws.Range('E20').Value = '  +4.41%  '
$ws.Range('D21').Value = '0.0₃0918'
$ws.Range('E21').Value = '  +0.89%  '
$ws.Range('E22').Value = '  +1.12%  '
$ws.Range('D23').Value = '''68.15'
$ws.Range('E23').Value = '  +1.25%  '
$ws.Range('D24').Value = '''243.83'
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').Value = '  +0.82%  '
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('D28').Value = '''24.13'
$ws.Range('E28').Value = '  -0.60%  '
$ws.Range('D29').Value = '''36.31'
$ws.Range('E29').Value = '  +6.56%  '
$ws.Range('E30').Value = '  +1.03%  '
$ws.Range('E31').Value = '  +1.46%  '
$ws.Range('D32').Value = '''161.58'
$ws.Range('E32').Value = '  +1.82%  '
$ws.Range('E33').Value = '  +3.25%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('E35').Value = '  +0.55%  '
$ws.Range('E36').Value = '  +1.30%  '
$ws.Range('E37').Value = '  +3.32%  '
$ws.Range('E38').Value = '  +2.38%  '
$ws.Range('E39').Value = '  -0.31%  '
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('E42').Value = '  +6.54%  '
$ws.Range('D43').Value = '2.020.48'
$ws.Range('E43').Value = '  -2.52%  '
$ws.Range('D44').Value = '''19.74'
$ws.Range('E44').Value = '  -0.43%  '
$ws.Range('E45').Value = '  +11.46%  '
$ws.Range('E46').Value = '  +1.81%  '
$ws.Range('D47').Value = '''10.23'
$ws.Range('E47').Value = '  -1.59%  '
$ws.Range('D48').Value = '''2.96'
$ws.Range('E48').Value = '  +1.58%  '
$ws.Range('D49').Value = '''53.54'
$ws.Range('E49').Value = '  +3.43%  '
$ws.Range('E50').Value = '  -0.68%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = '''1.15'
$ws.Range('E51').Value = '  -0.55%  '
